$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts all existing columns one to the right)
$ws.Columns("A:A").Insert()

# New column A header: "Match ID"
$ws.Range("A2").Value = "Match ID"

# Fill Match ID value (8) for the data rows
$ws.Range("A4:A19").Value = 8

# Bold the new "Match ID" column header + data (matches style used by column B's header/ID style)
$ws.Range("A2:A19").Font.Bold = $true

# Totals row keeps the Match ID value but without the bold style, and gets hidden
$ws.Range("A20").Value = 8
$ws.Rows(20).Hidden = $true

# Final selection highlighted by the author
$ws.Range("A2:A19").Select()
